# Updated audit code for modify & suspended
# - Rename the worksheet to reflect its actual contents (Service -> CW Roles mapping)
# - Move the active selection to reflect where the author was last working (C25)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "IDAM Mapping" sheet to "Service to CW Roles Mapping"
$ws.Name = "Service to CW Roles Mapping"

# Make sure it's the active sheet, then move/select the active cell to C25
$ws.Activate()
$ws.Range("C25").Select()
